$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new rows at the correct positions (formatting is inherited automatically)
$ws.Rows("6:7").Insert()
$ws.Rows("9:10").Insert()
$ws.Rows("15:16").Insert()

# Set cell values for all rows (header row 1 left untouched)
$ws.Range("A2").Value = "Instructions: In this set of questions we would like to know what the school reopening plans for all your school-aged children."
$ws.Range("C2").Value = "Developed by RAPID Team"
$ws.Range("D2").Value = "Current`n21, 23"

$ws.Range("A3").Value = "Are there children in your household who are in Kindergarten through 12th grade?"
$ws.Range("B3").Value = "• Yes`n• No"
$ws.Range("C3").Value = "Developed by RAPID Team"
$ws.Range("D3").Value = "Current`n21, 23"

$ws.Range("A4").Value = "[If yes to children in K-12]`nWhat has the school/school district your child(ren) plans to attend/currently attending decided regarding the school-year?"
$ws.Range("B4").Value = "• Open completely`n• Move to online only`n• A combination of in-person or some days and online on other days`n• Ability for parent to choose between online or in-person instruction`n• Other`n• Not applicable (e.g., homeschooled)`n• Don't know yet"
$ws.Range("C4").Value = "Developed by RAPID Team"
$ws.Range("D4").Value = "Current`n21, 23"

$ws.Range("A5").Value = "[If yes to children in K-12]`nFor the time your child will be learning remotely, who will be assisting with your child's online learning? Select all that apply."
$ws.Range("B5").Value = "• You/other parent/step-parent`n• Sibling 15 years or older`n• Grandparent`n• Other relative`n• Friend of parent `n• Other neighbor`n• Adult responsible for a group (e.g., pod, bubble)`n• Babysitter/nanny/au pair`n• Other`n• No one is able to do this`n• Not applicable"
$ws.Range("C5").Value = "Developed by RAPID Team"
$ws.Range("D5").Value = "Current`n21, 23"

$ws.Range("A6").Value = "[If yes to children in K-12]`nIf  your child has already begun the 2020-2021 school year, how are they reacting to online learning?"
$ws.Range("B6").Value = "•`tVery well`n•`tSomewhat well`n•`tWell`n•`tNot well `n•`tVery poorly"
$ws.Range("C6").Value = "Developed by RAPID Team"
$ws.Range("D6").Value = "Current 23"

$ws.Range("A7").Value = "[If yes to children in K-12]`nIncluding hours spent during weekdays and weekends, about how many hours did you spend on teaching activities with your school-aged child(ren) in this household during the last 7 days? Enter the total number of hours. If none, enter 0."
$ws.Range("B7").Value = "Open Response"
$ws.Range("C7").Value = "RAPID Modified `nU.S. Census Bureau Household Pulse Survey "
$ws.Range("D7").Value = "Current 23"

$ws.Range("A8").Value = "[If yes to children in K-12]`nWill you use a child care provider to help you when your child(ren) is not in school?"
$ws.Range("B8").Value = "• Yes`n• No`n• Maybe"
$ws.Range("C8").Value = "Developed by RAPID Team"
$ws.Range("D8").Value = "Current`n21, 23"

$ws.Range("A9").Value = "If you have a child that was due to be entering kindergarten this fall, have you decided to wait until next fall (2021) instead because of the pandemic?"
$ws.Range("B9").Value = "•`tYes`n•`tNo `n•`tNot applicable"
$ws.Range("C9").Value = "Developed by RAPID Team"
$ws.Range("D9").Value = "Current 23"

$ws.Range("A10").Value = "If yes, why? Select all that apply. "
$ws.Range("B10").Value = "•`tSafety`n•`tUncertain about the plan for school (in person/online)`n•`tNot able to manage online instruction for my child along with my other responsibilities (work, etc.)"
$ws.Range("C10").Value = "Developed by RAPID Team"
$ws.Range("D10").Value = "Current 23"

$ws.Range("A11").Value = "Does this affect the child care arrangements you have for your child(ren) age 5 and under?"
$ws.Range("B11").Value = "• Yes`n• No`n• Maybe`n• Not applicable"
$ws.Range("C11").Value = "Developed by RAPID Team"
$ws.Range("D11").Value = "Current`n21, 23"

$ws.Range("A12").Value = "Instructions: The following questions are about your child(ren)'s 0-5 learning."
$ws.Range("D12").Value = "Current`n21, 23"

$ws.Range("A13").Value = "For your child(ren) between the age of 0-5, which of the following are you doing to support their learning? Select all that apply."
$ws.Range("B13").Value = "• Using in-person games and activities at your home`n• Attending in-person learning outside of your home (daycare, childcare, etc.)`n• Using online resources (e.g., educational apps, educational TV shows, etc.)`n• Attending online classes/activities (facilitated by someone outside of your home, e.g., remote preschool, only story time, etc.)`n• Other`n• None of the above"
$ws.Range("C13").Value = "Developed by RAPID Team"
$ws.Range("D13").Value = "Current`n21, 23"

$ws.Range("A14").Value = "Who will be assisting with your child(s) 0-5 learning? Select all that apply."
$ws.Range("B14").Value = "• You/other parent/step-parent`n• Sibling 15 years or older`n• Grandparent`n• Other relative`n• Friend of parent `n• Other neighbor`n• Adult responsible for a group (e.g., pod, bubble)`n• Babysitter/nanny/au pair`n• Other`n• No one is able to do this`n• Not applicable"
$ws.Range("C14").Value = "Developed by RAPID Team"
$ws.Range("D14").Value = "Current`n21, 23"

$ws.Range("A15").Value = "Including hours spent during weekdays and weekends, about how many hours did you spend on face-to-face activities with your children 0-5 in this household during the last 7 days? Enter the total number of hours. If none, enter 0."
$ws.Range("B15").Value = "Open Response"
$ws.Range("C15").Value = "RAPID Modified `nU.S. Census Bureau Household Pulse Survey "
$ws.Range("D15").Value = "Current 23"

$ws.Range("A16").Value = "Including hours spent during weekdays and weekends, about how many hours did you spend working during the last 7 days? Enter the total number of hours. If none, enter 0 or NA if not working currently. "
$ws.Range("B16").Value = "Open Response"
$ws.Range("C16").Value = "RAPID Modified `nU.S. Census Bureau Household Pulse Survey "
$ws.Range("D16").Value = "Current 23"

# Fix up the C-column cells that need the wrapped/top-aligned style (copy format from col A of same row)
$ws.Range("A7").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Set row heights to match the target layout
$ws.Rows(2).RowHeight = 34
$ws.Rows(3).RowHeight = 34
$ws.Rows(4).RowHeight = 204
$ws.Rows(5).RowHeight = 221
$ws.Rows(6).RowHeight = 85
$ws.Rows(7).RowHeight = 68
$ws.Rows(8).RowHeight = 51
$ws.Rows(9).RowHeight = 51
$ws.Rows(10).RowHeight = 119
$ws.Rows(11).RowHeight = 68
$ws.Rows(12).RowHeight = 34
$ws.Rows(13).RowHeight = 255
$ws.Rows(14).RowHeight = 221
$ws.Rows(15).RowHeight = 51
$ws.Rows(16).RowHeight = 51

# Restore the selection/view state
$ws.Range("B15").Select()

